$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.808.87"
$ws.Range("E2").Value = "  +8.77%  "
$ws.Range("D3").Value = "1.952.09"
$ws.Range("E3").Value = "  +7.35%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'342.30"
$ws.Range("E5").Value = "  +3.00%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'0.4779"
$ws.Range("E7").Value = "  +4.75%  "
$ws.Range("D8").Value = "'0.4154"
$ws.Range("E8").Value = "  +9.11%  "
$ws.Range("D9").Value = "'48.50"
$ws.Range("E9").Value = "  +5.72%  "
$ws.Range("D10").Value = "'0.08270"
$ws.Range("E10").Value = "  +5.56%  "
$ws.Range("E11").Value = "  +9.25%  "
$ws.Range("D12").Value = "'22.73"
$ws.Range("E12").Value = "  +8.65%  "
$ws.Range("D13").Value = "'6.209"
$ws.Range("E13").Value = "  +6.91%  "
$ws.Range("D14").Value = "1.936.14"
$ws.Range("E14").Value = "  +5.35%  "
$ws.Range("D15").Value = "'7.433"
$ws.Range("E15").Value = "  +5.56%  "
$ws.Range("D16").Value = "'92.36"
$ws.Range("E16").Value = "  +3.65%  "
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "'0.00001065"
$ws.Range("E18").Value = "  +4.86%  "
$ws.Range("D19").Value = "'0.06675"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "'18.13"
$ws.Range("E20").Value = "  +6.44%  "
$ws.Range("D21").Value = "'0.9998"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").Value = "29.764.61"
$ws.Range("E22").Value = "  +8.65%  "
$ws.Range("D23").Value = "'5.603"
$ws.Range("E23").Value = "  +6.28%  "
$ws.Range("E24").Value = "  +4.75%  "
$ws.Range("D25").Value = "'2.280"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").Value = "2.169.75"
$ws.Range("E26").Value = "  +5.61%  "
$ws.Range("D27").Value = "'160.85"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("D28").Value = "'20.20"
$ws.Range("E28").Value = "  +4.97%  "
$ws.Range("D29").Value = "'2.192"
$ws.Range("E29").Value = "  +7.88%  "
$ws.Range("D30").Value = "'5.706"
$ws.Range("E30").Value = "  +8.78%  "
$ws.Range("D31").Value = "'122.68"
$ws.Range("E31").Value = "  +4.49%  "
$ws.Range("D32").Value = "'1.021"
$ws.Range("E32").Value = "  +10.02%  "
$ws.Range("D33").Value = "'0.09643"
$ws.Range("E33").Value = "  +3.60%  "
$ws.Range("D34").Value = "'1.482"
$ws.Range("E34").Value = "  +13.08%  "
$ws.Range("D35").Value = "'3.680"
$ws.Range("E35").Value = "  +3.32%  "
$ws.Range("D36").Value = "'5.528"
$ws.Range("E36").Value = "  +6.46%  "
$ws.Range("D37").Value = "'0.06315"
$ws.Range("E37").Value = "  +7.28%  "
$ws.Range("D38").Value = "'0.02337"
$ws.Range("E38").Value = "  +7.50%  "
$ws.Range("D39").Value = "'8.558"
$ws.Range("E39").Value = "  +6.16%  "
$ws.Range("D40").Value = "'1.200"
$ws.Range("E40").Value = "  +5.51%  "
$ws.Range("D41").Value = "'0.6125"
$ws.Range("E41").Value = "  +7.02%  "
$ws.Range("D42").Value = "'10.74"
$ws.Range("E42").Value = "  +8.85%  "
$ws.Range("D43").Value = "'0.1900"
$ws.Range("E43").Value = "  +5.09%  "
$ws.Range("D44").Value = "'0.9997"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "'1.291"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("D46").Value = "'2.397"
$ws.Range("D47").Value = "'12.60"
$ws.Range("E47").Value = "  +7.31%  "
$ws.Range("D48").Value = "'0.5734"
$ws.Range("E48").Value = "  +6.69%  "
$ws.Range("D49").Value = "'2.002"
$ws.Range("E49").Value = "  +7.30%  "
$ws.Range("E50").Value = "  +12.94%  "
$ws.Range("D51").Value = "'114.58"
$ws.Range("E51").Value = "  +4.13%  "
